# Rotate the data rows 2-5 of the "Artfynd" sheet:
#   new row2 <- old row3
#   new row3 <- old row4
#   new row4 <- old row5
#   new row5 <- old row2
# (row numbers stay put; only the record payload moves, one position "up",
#  with the top record wrapping around to the bottom.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 5
$lastCol  = 51   # column AY

# Regex used to detect "yyyy-mm-dd" text so we can stop Excel's COM layer
# from auto-coercing it into a date serial number when we write it back.
$datePattern = '^\d{4}-\d{2}-\d{2}$'

# 1) Snapshot every cell in rows 2..5 (raw .Value2, unaffected by display
#    formatting) BEFORE any writes happen, since the rotation would
#    otherwise clobber source data before it's been read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Destination row -> source row mapping implementing the rotation.
$mapping = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($destRow in 2, 3, 4, 5) {
    $srcRow = $mapping[$destRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $snapshot["$srcRow,$c"]
        $cell = $ws.Cells.Item($destRow, $c)

        if ($null -eq $val -or $val -eq "") {
            # Source cell was blank -> make sure destination ends up blank too.
            $cell.Value = ""
        }
        elseif ($val -is [string] -and $val -match $datePattern) {
            # Force text so Excel doesn't reinterpret "2023-03-11" etc. as a
            # date serial number; restore default styling afterwards so we
            # don't leave a stray number-format behind.
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        }
        else {
            $cell.Value = $val
        }
    }
}
